$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-23
# from 46074 (2026-02-21) to 46075 (2026-02-22)
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 46075
}
